# Webshop Checkliste: add "Punkte als Rabatt einlösen" checklist item (Kasse/Punkte)
# at row 46, shifting all subsequent checklist rows down by one. Mark several newly
# reviewed rows as done ("Ja") in column E (Datenschutzprüfung / Warenkorbrabatt bei
# 5/10 Artikeln). Update current selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 46; existing rows 46-60 shift down to 47-61.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new checklist entry.
$ws.Cells.Item(46, 1).Value = 10
$ws.Cells.Item(46, 2).Value = "Kasse"
$ws.Cells.Item(46, 3).Value = "Punkte"
$ws.Cells.Item(46, 4).Value = "Punkte als Rabatt einlösen"
$ws.Cells.Item(46, 5).Value = "Ja"
$ws.Cells.Item(46, 6).Value = "KIMI"

# Mark additional rows as completed ("Ja") in column E: 5%/10% cart discount checks
# (rows 42-43, unaffected by the insert) and the three Kasse rows now shifted to 47-49
# (shipping, cost, and privacy-checkbox / Datenschutzprüfung checks).
$ws.Cells.Item(42, 5).Value = "Ja"
$ws.Cells.Item(43, 5).Value = "Ja"
$ws.Cells.Item(47, 5).Value = "Ja"
$ws.Cells.Item(48, 5).Value = "Ja"
$ws.Cells.Item(49, 5).Value = "Ja"

# Update the recorded selection to match the author's final cursor position.
$ws.Range("E50").Select()
